$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NOTES")
Write-Output ($ws | Get-Member | Out-String)
